$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'302.13"
$ws.Range("E2").Value = "'-4.01%"

$ws.Range("D3").Value = "'35.21"
$ws.Range("E3").Value = "'-2.10%"

$ws.Range("D4").Value = "'5.043"
$ws.Range("E4").Value = "'-1.60%"

$ws.Range("D5").Value = "'0.07984"
$ws.Range("E5").Value = "'-1.72%"

$ws.Range("D6").Value = "'1.936"
$ws.Range("E6").Value = "'-8.74%"

$ws.Range("B7").Value = 'KuCoinToken'
$ws.Range("C7").Value = 'https://coinranking.com/coin/LOO6LmXd7G84Z+kucointoken-kcs'
$ws.Range("D7").Value = "'7.741"
$ws.Range("E7").Value = "'-3.50%"

$ws.Range("B8").Value = 'BTSEToken'
$ws.Range("C8").Value = 'https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse'
$ws.Range("D8").Value = "'2.911"
$ws.Range("E8").Value = "'7.64%"

$ws.Range("B9").Value = 'MXToken'
$ws.Range("C9").Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range("D9").Value = "'0.9238"
$ws.Range("E9").Value = "'-0.64%"

$ws.Range("B10").Value = 'LiechtensteinCryptoassetsExchange'
$ws.Range("C10").Value = 'https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx'
$ws.Range("D10").Value = "'0.1296"
$ws.Range("E10").Value = "'27.27%"

$ws.Range("B11").Value = 'WazirX'
$ws.Range("C11").Value = 'https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx'
$ws.Range("D11").Value = "'0.1846"
$ws.Range("E11").Value = "'-1.43%"

$ws.Range("B12").Value = 'MandalaExchangeToken'
$ws.Range("C12").Value = 'https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx'
$ws.Range("D12").Value = "'0.09628"
$ws.Range("E12").Value = "'5.11%"

$ws.Range("B13").Value = 'BitrueCoin'
$ws.Range("C13").Value = 'https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr'
$ws.Range("D13").Value = "'0.03633"
$ws.Range("E13").Value = "'1.30%"

$ws.Range("B14").Value = 'BitMartToken'
$ws.Range("C14").Value = 'https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx'
$ws.Range("D14").Value = "'0.09851"
$ws.Range("E14").Value = "'-0.71%"

$ws.Range("B15").Value = 'BitForexToken'
$ws.Range("C15").Value = 'https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf'
$ws.Range("D15").Value = "'0.001395"
$ws.Range("E15").Value = "'-3.07%"

$ws.Range("B16").Value = 'TigerCash'
$ws.Range("C16").Value = 'https://coinranking.com/coin/6hIn06L2+tigercash-tch'
$ws.Range("D16").Value = "'0.005821"
$ws.Range("E16").Value = "'1.05%"

$ws.Range("B17").Value = 'LEO'
$ws.Range("C17").Value = 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'
$ws.Range("D17").Value = "'3.505"
$ws.Range("E17").Value = "'0.92%"

$ws.Range("B18").Value = 'GateToken'
$ws.Range("C18").Value = 'https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt'
$ws.Range("D18").Value = "'4.047"
$ws.Range("E18").Value = "'-2.37%"

$ws.Range("D19").Value = "'0.3430"
$ws.Range("E19").Value = "'1.79%"

$ws.Range("D20").Value = "'0.1311"
$ws.Range("E20").Value = "'-1.74%"

$ws.Range("D21").Value = "'5.046"
$ws.Range("E21").Value = "'-1.73%"

$ws.Range("D22").Value = "'0.2401"
$ws.Range("E22").Value = "'9.39%"

$ws.Range("D23").Value = "'0.04525"
$ws.Range("E23").Value = "'-1.16%"

$ws.Range("D24").Value = "'0.001218"
$ws.Range("E24").Value = "'-2.64%"

$ws.Range("D25").Value = "'0.004815"
$ws.Range("E25").Value = "'2.00%"

$ws.Range("D26").Value = "'0.0001251"
$ws.Range("E26").Value = "'-0.28%"

$ws.Range("E27").Value = "'-33.45%"

$ws.Range("D39").Value = "'0.01900"
$ws.Range("E39").Value = "'-4.20%"

$ws.Range("D40").Value = "'0.04699"
$ws.Range("E40").Value = "'-4.32%"

$ws.Range("D41").Value = "'0.007554"
$ws.Range("E41").Value = "'-4.08%"

$ws.Range("D42").Value = "'0.009670"
$ws.Range("E42").Value = "'23.47%"

$ws.Range("E43").Value = "'-5.20%"

$ws.Range("D44").Value = "'0.002112"
$ws.Range("E44").Value = "'0.05%"

$ws.Range("D45").Value = "'0.01088"
$ws.Range("E45").Value = "'-6.50%"

$ws.Range("D46").Value = "'0.00006234"

$ws.Range("E47").Value = "'-0.24%"

$ws.Range("E48").Value = "'65.22%"

$ws.Range("E49").Value = "'-21.89%"

$ws.Range("D50").Value = "'0.00002102"
$ws.Range("E50").Value = "'-0.24%"

$ws.Range("E51").Value = "'-0.24%"
